# Insert the new "Git and GitHub" topic section (one bold heading paragraph
# plus five ListParagraph sub-bullets) right after the existing
# "Reading and writing to text and CSV files" bullet (the last item under
# "File Handling"), mirroring the rest of the outline's structure/numbering.

$d = $word.ActiveDocument

# Locate the anchor paragraph by scanning for its exact text (robust across
# hosts, avoids relying on stateful Find.Found after re-deriving .Content).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Reading and writing to text and CSV files") {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "Anchor paragraph 'Reading and writing to text and CSV files' not found"
}

# Create a fresh empty paragraph right after the anchor, then replace its
# contents with the exact OOXML for the six new paragraphs (this preserves
# run-level formatting, proofErr spell-check hints, and literal punctuation
# exactly as authored).
$anchorRange = $d.Paragraphs($anchorIndex).Range
$anchorRange.InsertParagraphAfter()

$newPara = $d.Paragraphs($anchorIndex + 1)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>Git</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> and GitHub</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>How Git and GitHub work</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Explanation of the process for including changes in a local repository</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>(working directory, staging area, and local repository)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Commands to include changes in a local repository </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">git status, git add, and git commit -m </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>‘</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>comment</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>’</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>How to push a local repository to a remote repository</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>¿</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>What is the library called to integrate Git and Python?</w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newPara.Range.InsertXML($xml)

Write-Output "Inserted Git and GitHub section after paragraph $anchorIndex"
